# Generate Report for Handback
# - Overview sheet: update status for the 165a8ba4... row from
#   "Ready for handoff" to "Handback transform failed"
# - zh-cn / de-de sheets: record the Error Detail for that row explaining
#   why the handback transform failed (filename mismatch).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The "Ready for handoff" status for file 165a8ba4-24a5-4680-bf1d-72ec2568d889
# is shared text used on the Overview sheet (B3 = zh-cn status, C3 = de-de
# status) and on the per-locale Status column (C3) of the zh-cn and de-de
# sheets. The handback transform failed for this file, so every occurrence
# becomes "Handback transform failed".
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: add Error Detail (column L) on row 3 explaining the failure
$zhcn.Range("L3").Value = "Handback file name: dmza43vp.ila is different with handoff file name: 165a8ba4-24a5-4680-bf1d-72ec2568d889.c165ea14b96dd7463b41e31dfc6dcc9d007f74c0.zh-cn."

# de-de sheet: add Error Detail (column L) on row 3 explaining the failure
$dede.Range("L3").Value = "Handback file name: dmza43vp.ila is different with handoff file name: 165a8ba4-24a5-4680-bf1d-72ec2568d889.c165ea14b96dd7463b41e31dfc6dcc9d007f74c0.de-de."
